$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking values
# (e.g. "1.00", "67.13", "0.0520") are preserved exactly as strings
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '34.470.17'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = '1.800.80'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').Value = '224.27'
$ws.Range('D6').Value = '0.598'
$ws.Range('E6').Value = '  +0.33%  '
$ws.Range('E7').Value = '  +0.32%  '
$ws.Range('D8').Value = '40.99'
$ws.Range('E8').Value = '  +13.38%  '
$ws.Range('D9').Value = '0.291'
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  -0.88%  '
$ws.Range('D11').Value = '0.0997'
$ws.Range('E11').Value = '  +3.53%  '
$ws.Range('D12').Value = '2.060.43'
$ws.Range('E12').Value = '  +0.35%  '
$ws.Range('D13').Value = '1.793.81'
$ws.Range('E13').Value = '  +0.14%  '
$ws.Range('D14').Value = '10.92'
$ws.Range('E14').Value = '  -2.23%  '
$ws.Range('D15').Value = '34.435.86'
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('D16').Value = '0.626'
$ws.Range('E16').Value = '  -0.54%  '
$ws.Range('D17').Value = '4.40'
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').Value = '67.13'
$ws.Range('E18').Value = '  -1.93%  '
$ws.Range('D19').Value = '239.63'
$ws.Range('E19').Value = '  -1.07%  '
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('D21').Value = '11.11'
$ws.Range('E21').Value = '  -1.11%  '
$ws.Range('E22').Value = '  +0.29%  '
$ws.Range('D23').Value = '4.16'
$ws.Range('E23').Value = '  +2.55%  '
$ws.Range('D24').Value = '2.16'
$ws.Range('E24').Value = '  -1.15%  '
$ws.Range('D25').Value = '172.10'
$ws.Range('E25').Value = '  +0.77%  '
$ws.Range('D26').Value = '7.65'
$ws.Range('E26').Value = '  -2.62%  '
$ws.Range('D27').Value = '17.35'
$ws.Range('E27').Value = '  +0.57%  '
$ws.Range('E28').Value = '  -0.19%  '
$ws.Range('E29').Value = '  +0.27%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = '3.77'
$ws.Range('E30').Value = '  +0.25%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '1.22'
$ws.Range('E31').Value = '  +0.09%  '
$ws.Range('D32').Value = '3.85'
$ws.Range('E32').Value = '  -0.36%  '
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('D34').Value = '1.77'
$ws.Range('E34').Value = '  +1.00%  '
$ws.Range('D35').Value = '1.321.00'
$ws.Range('E35').Value = '  -2.67%  '
$ws.Range('D36').Value = '0.645'
$ws.Range('E36').Value = '  +0.35%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = '1.06'
$ws.Range('E37').Value = '  +0.71%  '
$ws.Range('B38').Value = 'Aave'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D38').Value = '85.75'
$ws.Range('E38').Value = '  +6.73%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '2.36'
$ws.Range('E39').Value = '  +0.55%  '
$ws.Range('D40').Value = '14.79'
$ws.Range('E40').Value = '  +13.54%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '0.0187'
$ws.Range('E41').Value = '  +1.54%  '
$ws.Range('D42').Value = '1.23'
$ws.Range('E42').Value = '  +6.25%  '
$ws.Range('E43').Value = '  +0.76%  '
$ws.Range('E44').Value = '  +0.27%  '
$ws.Range('E45').Value = '  +0.55%  '
$ws.Range('D46').Value = '0.0520'
$ws.Range('E46').Value = '  +5.11%  '
$ws.Range('D47').Value = '1.961.15'
$ws.Range('E47').Value = '  +0.35%  '
$ws.Range('D48').Value = '5.80'
$ws.Range('E48').Value = '  +1.53%  '
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').Value = '  +0.32%  '
$ws.Range('D50').Value = '100.22'
$ws.Range('E50').Value = '  -1.49%  '
$ws.Range('D51').Value = '0.0609'
$ws.Range('E51').Value = '  +1.16%  '

# Restore the original (default) cell style for column D now that
# the text values are set, so no stray number-format style remains
# applied to the data cells.
$ws.Range("D2:D51").Style = "Normal"
